$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# GAME_ID changes from 0022000535 to 0022000098 for both rows.
# Force Text format first so the leading zeros in the numeric-looking
# string are preserved instead of Excel auto-converting it to a number.
$ws.Range("B2:B3").NumberFormat = "@"
$ws.Range("B2").Value = "0022000098"
$ws.Range("B3").Value = "0022000098"

# Row 2: now the Houston Rockets (TEAM_ID 1610612745)
$ws.Range("C2").Value = 1610612745
$ws.Range("D2").Value = "Rockets"
$ws.Range("E2").Value = "HOU"
$ws.Range("F2").Value = "Houston"
$ws.Range("H2").Value = 0.461
$ws.Range("I2").Value = 0.468
$ws.Range("J2").Value = 0.123
$ws.Range("K2").Value = 0.137
$ws.Range("L2").Value = 0.573
$ws.Range("M2").Value = 0.244
$ws.Range("N2").Value = 0.182
$ws.Range("O2").Value = 0.326

# Row 3: now the Dallas Mavericks (TEAM_ID 1610612742)
$ws.Range("C3").Value = 1610612742
$ws.Range("D3").Value = "Mavericks"
$ws.Range("E3").Value = "DAL"
$ws.Range("F3").Value = "Dallas"
$ws.Range("H3").Value = 0.573
$ws.Range("I3").Value = 0.244
$ws.Range("J3").Value = 0.182
$ws.Range("K3").Value = 0.233
$ws.Range("L3").Value = 0.461
$ws.Range("M3").Value = 0.468
$ws.Range("N3").Value = 0.123
$ws.Range("O3").Value = 0.176
